$d = $word.ActiveDocument
$CR = [char]13

# We edit from the end of the document towards the beginning so that
# paragraph indices of not-yet-processed paragraphs remain stable.

# ---------------------------------------------------------------------
# Paragraph 18 (empty) -> "Extract all code and input images..."
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "Extract all code and input images in this .zip file, and run main.m. To save time, the images have been presampled. To resample images, edit line 10 of main.m to read resample=1. Resampling takes several minutes."

# ---------------------------------------------------------------------
# Paragraph 17: "Conclusions" -> "Instructions for Running Code"
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)
$r17 = $d.Range($p17.Range.Start, $p17.Range.End - 1)
$r17.Text = "Instructions for Running Code"

# ---------------------------------------------------------------------
# Paragraph 16 (empty) -> 3 new paragraphs:
#   "Over the course of this project..."
#   "The algorithms run very quickly..."
#   "Of course, some of the overhead..." (contains a lastRenderedPageBreak
#    right before the final run "network. I had intended...")
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$p16.Range.Text = "Over the course of this project, I found that both of these classifying algorithms would be suitable for classifying bananas. Even with images from disparate sources with varying resolution and color cast, the results were acceptable. With more tuning and properly calibrated input images, either of these solutions would be suitable for industrial use." + $CR + "The algorithms run very quickly, but the preprocessing of images (eliminating background pixels, sampling, etc) runs very slowly, because it is written inefficiently. Since that part of the code was not the focus of the project, little development time was applied to making it robust and efficient. If that part of the code were optimized, it would be reasonable for either of these methods to train from an extensive input data set in under a minute and to classify an input in under one second." + $CR
$p18b = $d.Paragraphs.Item(18)
$p18b.Range.Text = "Of course, some of the overhead of sampling could be eliminated if the pixel grouping process were unnecessary. Ideally, the solution would be to use a perceptron-based neural network for this classification task, because it would perform the grouping itself as part of the architecture of the network. I had intended to write a perceptron classifier as a third method for comparison, but time did not allow me to complete it."

# Insert the lastRenderedPageBreak marker right before "network. I had intended"
$p18c = $d.Paragraphs.Item(18)
$breakPoint = $p18c.Range.Start + ("Of course, some of the overhead of sampling could be eliminated if the pixel grouping process were unnecessary. Ideally, the solution would be to use a perceptron-based neural network for this classification task, because it would perform the grouping itself as part of the architecture of the ").Length
$rb = $d.Range($breakPoint, $breakPoint)
$ooxmlBreak = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rb.InsertXML($ooxmlBreak)

# ---------------------------------------------------------------------
# Paragraph 15: "Perceptron Classification of Random Pixel Samples"
#   -> "Finding" + "s" (two runs, heading style)
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$r15 = $d.Range($p15.Range.Start, $p15.Range.End - 1)
$r15.Text = "Findings"

# ---------------------------------------------------------------------
# Paragraph 14 (empty) -> 2 new paragraphs:
#  "I also applied K-Means clustering..."
#  "With 20 clusters..."
# ---------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "I also applied K-Means clustering to this problem, and had similar success. The input to the k-means algorithm was the same set of pixel-group-count vectors used as an input to the self-ordered mapping algorithm in my second attempt. The K-means clustering algorithm was then run for (up to) 1000 iterations (it converged much faster in every test)." + $CR
$p15b = $d.Paragraphs.Item(15)
$p15b.Range.Text = "With 20 clusters, the k-means algorithm was able to correctly identify 8 out of 9 ripe bananas, misclassifying one as under ripe. It performed worst on under ripe bananas, classifying 3 out of 5 correctly and misclassifying two as ripe. It performed the best on the overripe bananas, correctly classifying 4 out of 5 and only failing on the greenish-tinted banana that the self-ordered mapping algorithm failed on."

# ---------------------------------------------------------------------
# Paragraph 13: "K-Means Clustering of Pixel Bin Counts"
#   remove the lastRenderedPageBreak marker (text unchanged)
# ---------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$r13 = $d.Range($p13.Range.Start, $p13.Range.End - 1)
$r13.Text = "TEMP_PLACEHOLDER"
$p13b = $d.Paragraphs.Item(13)
$r13b = $d.Range($p13b.Range.Start, $p13b.Range.End - 1)
$r13b.Text = "K-Means Clustering of Pixel Bin Counts"

# ---------------------------------------------------------------------
# Paragraph 12 (empty) -> 2 new paragraphs:
#  "In order to solve this problem..."
#  "With this input, the results..."
# ---------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "In order to solve this problem, I changed the input vector to the algorithm. After randomly sampling the input images, the pixels were divided into groups by their hue value. The total number of pixels in each group was counted, and the list of counts was used as the input vector for each pattern. So for each image, a vector of pixel counts was created and fed into the self-ordered mapping network." + $CR
$p13c = $d.Paragraphs.Item(13)
$p13c.Range.Text = "With this input, the results were much more promising. Using 100 pixel groups and 5 output bins, the algorithm can correctly identify 7 out of 9 ripe bananas, misclassifying 2 as overripe. It correctly classifies 3 out of 5 under ripe bananas, misclassifying 2 as ripe. It correctly classifies 4 out of 5 overripe bananas, and the one it misclassifies has a greenish tint to the image, which caused all tests to classify it as under ripe. In a real situation, controlled lighting and camera calibration would eliminate this problem."

# ---------------------------------------------------------------------
# Paragraph 11: "Self-Ordered Mapping of Pixel Bin Counts"
#   -> "Self-Ordered Mapping of Pixel " + "Group" + " Counts"
# ---------------------------------------------------------------------
$p11 = $d.Paragraphs.Item(11)
$r11 = $d.Range($p11.Range.Start, $p11.Range.End - 1)
$r11.Text = "Self-Ordered Mapping of Pixel Group Counts"

# ---------------------------------------------------------------------
# Paragraph 10: "Since the pixels were from random locations..."
#   add a lastRenderedPageBreak marker at the very start of the run
# ---------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$rb10 = $d.Range($p10.Range.Start, $p10.Range.Start)
$rb10.InsertXML($ooxmlBreak)

# ---------------------------------------------------------------------
# Paragraph 6: "Given a set of banana images, several thousand..."
#   "several" (highlighted) -> "one" (no highlight)
# ---------------------------------------------------------------------
$find6 = $d.Content.Find
$find6.ClearFormatting()
$found6 = $find6.Execute("several", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found6) {
    $find6.Parent.HighlightColorIndex = 0
    $find6.Parent.Text = "one"
}

# Insert the new validation-set paragraph right after paragraph 6
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7new = $d.Paragraphs.Item(7)
$p7new.Range.Text = "The set of images were sampled one more time to generate a validation set. This set of 19 images (9 ripe, 5 under ripe, and 5 overripe) were used to evaluate the classification abilities of the network."
$p7new.Range.ParagraphFormat.Style = "Normal"

# ---------------------------------------------------------------------
# Paragraph 5: "Since isolating the banana(s)..."
#   prepend a new leading sentence as its own run
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$rstart5 = $d.Range($p5.Range.Start, $p5.Range.Start)
$rstart5.InsertBefore("The images used in this project were all obtained from the Internet, and contained various background scenery and other extraneous information. ")

# ---------------------------------------------------------------------
# Paragraph 3: "Several of the algorithms..." rewrite with split runs
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = "Several of the algorithms in this class can be applied to image processing tasks, and one such task is determination of fruit ripeness. In this report, I examine the results of applying Self-Ordered Mapping and K-Means Clustering to the challenge of classifying bananas as under ripe, ripe, or overripe. These methods, if refined, could be employed in automated fruit packing plants and in other parts of the food industry to automatically sort fruit."

# ---------------------------------------------------------------------
# Paragraph 1: insert "By Edward Venator" as new paragraph 2
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$p2new = $d.Paragraphs.Item(2)
$p2new.Range.Text = "By Edward Venator"
$p2new.Range.ParagraphFormat.Style = "Normal"

Write-Output "All edits applied"
